$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New variant/fitness dataset (row 2 keeps "4F" unchanged, rest replaced/reordered)
$data = @(
    @("4F", 1.0455),
    @("4H", 0.989),
    @("4L", 0.9845),
    @("4N", 1.0115),
    @("4Q", 0.994),
    @("4V", 0.99775),
    @("5M", 1.038),
    @("8C", 0.8985),
    @("8Q", 0.913),
    @("12C", 0.7385),
    @("12M", 1.012),
    @("71C", 1)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
